$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits at the
#    end of the "\dt - list the tables in the current database"
#    paragraph. It is a hidden bookmark (name starts with "_") so it
#    is not enumerated via Bookmarks.Count, but it can still be
#    reached (and removed) by name.
# ------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
    # Bookmark already absent - nothing to do.
}

# ------------------------------------------------------------------
# 2) Append the new "Linux Shell" block at the end of the document,
#    after the existing "psql practicedb" paragraph, finishing with
#    the "_GoBack" bookmark re-anchored on the final new paragraph.
# ------------------------------------------------------------------
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$newBlockXml = @"
<w:p $wNs/>
<w:p $wNs/>
<w:p $wNs/>
<w:p $wNs><w:r><w:t>Linux Shell</w:t></w:r></w:p>
<w:p $wNs/>
<w:p $wNs/>
<w:p $wNs><w:r><w:t>-</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>u :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> user</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>-</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>d :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> database</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>-c : command</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML($newBlockXml)
